$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values of rows 4, 5, 6 for the columns that change.
$cols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R", "Z", "AB")

$row4 = @{}
$row5 = @{}
$row6 = @{}

foreach ($col in $cols) {
    $row4[$col] = $ws.Range("${col}4").Value()
    $row5[$col] = $ws.Range("${col}5").Value()
    $row6[$col] = $ws.Range("${col}6").Value()
}

# Perform the cyclic rotation:
# new row4 = old row6
# new row5 = old row4
# new row6 = old row5
foreach ($col in $cols) {
    $ws.Range("${col}4").Value = $row6[$col]
    $ws.Range("${col}5").Value = $row4[$col]
    $ws.Range("${col}6").Value = $row5[$col]
}
